$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "89.947.36"
$ws.Range("E2").Value = "  +0.07%  "
$ws.Range("D3").Value = "3.111.07"
$ws.Range("E3").Value = "  +0.96%  "
$ws.Range("E4").Value = "  +0.31%  "
$ws.Range("D5").Value = "'235.49"
$ws.Range("E5").Value = "  +0.95%  "
$ws.Range("D6").Value = "'614.27"
$ws.Range("E6").Value = "  -0.64%  "
$ws.Range("D7").Value = "'1.08"
$ws.Range("E7").Value = "  +4.02%  "
$ws.Range("D8").Value = "'0.365"
$ws.Range("E8").Value = "  +2.32%  "
$ws.Range("E9").Value = "  -0.02%  "
$ws.Range("D10").Value = "3.114.63"
$ws.Range("E10").Value = "  +1.19%  "
$ws.Range("D11").Value = "'0.728"
$ws.Range("E11").Value = "  +2.24%  "
$ws.Range("E12").Value = "  +2.92%  "
$ws.Range("D13").Value = "'0.0000242"
$ws.Range("E13").Value = "  -0.88%  "
$ws.Range("D14").Value = "'34.63"
$ws.Range("E14").Value = "  -0.68%  "
$ws.Range("D15").Value = "'5.47"
$ws.Range("E15").Value = "  +2.50%  "
$ws.Range("D16").Value = "90.166.01"
$ws.Range("E16").Value = "  +0.39%  "
$ws.Range("D17").Value = "3.698.97"
$ws.Range("E17").Value = "  +1.24%  "
$ws.Range("D18").Value = "3.148.11"
$ws.Range("E18").Value = "  +2.07%  "
$ws.Range("D19").Value = "'3.63"
$ws.Range("E19").Value = "  -3.99%  "
$ws.Range("D20").Value = "'14.88"
$ws.Range("E20").Value = "  +8.18%  "
$ws.Range("D21").Value = "'5.76"
$ws.Range("E21").Value = "  +6.25%  "
$ws.Range("D22").Value = "'0.0000198"
$ws.Range("E22").Value = "  -4.92%  "
$ws.Range("D23").Value = "'435.55"
$ws.Range("E23").Value = "  +1.40%  "
$ws.Range("D24").Value = "'8.92"
$ws.Range("E24").Value = "  +2.22%  "
$ws.Range("D25").Value = "'5.67"
$ws.Range("E25").Value = "  +2.00%  "
$ws.Range("B26").Value = "Aptos"
$ws.Range("C26").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D26").Value = "'11.61"
$ws.Range("E26").Value = "  -0.81%  "
$ws.Range("B27").Value = "Litecoin"
$ws.Range("C27").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D27").Value = "'81.53"
$ws.Range("E27").Value = "  -5.15%  "
$ws.Range("D28").Value = "3.342.79"
$ws.Range("E28").Value = "  +1.90%  "
$ws.Range("D29").Value = "'0.998"
$ws.Range("E29").Value = "  -0.14%  "
$ws.Range("D30").Value = "'0.126"
$ws.Range("E30").Value = "  +46.18%  "
$ws.Range("D31").Value = "'0.226"
$ws.Range("E31").Value = "  +18.24%  "
$ws.Range("E32").Value = "  +7.91%  "
$ws.Range("D33").Value = "'9.16"
$ws.Range("E33").Value = "  +0.89%  "
$ws.Range("B34").Value = "Kaspa"
$ws.Range("C34").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D34").Value = "'0.168"
$ws.Range("E34").Value = "  +12.87%  "
$ws.Range("B35").Value = "Binance-PegBSC-USD"
$ws.Range("C35").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D35").Value = "'0.924"
$ws.Range("E35").Value = "  -8.03%  "
$ws.Range("D36").Value = "'7.60"
$ws.Range("E36").Value = "  +8.55%  "
$ws.Range("D37").Value = "'25.81"
$ws.Range("E37").Value = "  +1.34%  "
$ws.Range("B38").Value = "Bittensor"
$ws.Range("C38").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D38").Value = "'496.16"
$ws.Range("E38").Value = "  +0.15%  "
$ws.Range("B39").Value = "PancakeSwap"
$ws.Range("C39").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D39").Value = "'1.91"
$ws.Range("E39").Value = "  +1.69%  "
$ws.Range("E40").Value = "  +5.73%  "
$ws.Range("D41").Value = "'0.439"
$ws.Range("E41").Value = "  +10.68%  "
$ws.Range("D42").Value = "'3.72"
$ws.Range("E42").Value = "  +3.59%  "
$ws.Range("E43").Value = "  -8.23%  "
$ws.Range("D44").Value = "'22.07"
$ws.Range("E44").Value = "  -0.07%  "
$ws.Range("E45").Value = "  +0.00%  "
$ws.Range("D46").Value = "'158.34"
$ws.Range("E46").Value = "  +6.05%  "
$ws.Range("D47").Value = "'0.703"
$ws.Range("E47").Value = "  +4.38%  "
$ws.Range("E48").Value = "  +1.27%  "
$ws.Range("E49").Value = "  +3.67%  "
$ws.Range("D50").Value = "'43.82"
$ws.Range("E50").Value = "  -1.26%  "
$ws.Range("E51").Value = "  +0.62%  "
